# Update TPM values and sending-cluster labels for Apln-Aplnr LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.25093
$ws.Range("H2").Value = 60.75279
$ws.Range("I2").Value = 0.7748170638893737
$ws.Range("J2").Value = 0.7748170638893737
$ws.Range("M2").Value = 186.4969126666666
$ws.Range("N2").Value = 559.490738
$ws.Range("O2").Value = 0.9838307803020685
$ws.Range("P2").Value = 0.9838307803020684
$ws.Range("Q2").Value = 3776.73592362878
$ws.Range("R2").Value = 33990.62331265902
$ws.Range("S2").Value = 0.7622888765576401
$ws.Range("T2").Value = 0.7622888765576401
$ws.Range("G3").Value = 20.25093
$ws.Range("H3").Value = 60.75279
$ws.Range("I3").Value = 0.7748170638893737
$ws.Range("J3").Value = 0.7748170638893737
$ws.Range("O3").Value = 0.004224423017480373
$ws.Range("P3").Value = 0.004224423017480372
$ws.Range("Q3").Value = 16.2167422347
$ws.Range("R3").Value = 145.9506801123
$ws.Range("S3").Value = 0.00327315503903083
$ws.Range("T3").Value = 0.00327315503903083
$ws.Range("G4").Value = 20.25093
$ws.Range("H4").Value = 60.75279
$ws.Range("I4").Value = 0.7748170638893737
$ws.Range("J4").Value = 0.7748170638893737
$ws.Range("M4").Value = 1.711075333333334
$ws.Range("N4").Value = 5.133226000000001
$ws.Range("O4").Value = 0.009026468890440984
$ws.Range("P4").Value = 0.009026468890440982
$ws.Range("Q4").Value = 34.65086680006
$ws.Range("R4").Value = 311.8578012005401
$ws.Range("S4").Value = 0.006993862122980256
$ws.Range("T4").Value = 0.006993862122980254
$ws.Range("G5").Value = 20.25093
$ws.Range("H5").Value = 60.75279
$ws.Range("I5").Value = 0.7748170638893737
$ws.Range("J5").Value = 0.7748170638893737
$ws.Range("M5").Value = 0.1126546666666667
$ws.Range("N5").Value = 0.337964
$ws.Range("O5").Value = 0.00059428934788552
$ws.Range("P5").Value = 0.0005942893478855199
$ws.Range("Q5").Value = 2.28136176884
$ws.Range("R5").Value = 20.53225591956
$ws.Range("S5").Value = 0.0004604655276293891
$ws.Range("T5").Value = 0.0004604655276293891
$ws.Range("G6").Value = 20.25093
$ws.Range("H6").Value = 60.75279
$ws.Range("I6").Value = 0.7748170638893737
$ws.Range("J6").Value = 0.7748170638893737
$ws.Range("M6").Value = 0.4405493333333334
$ws.Range("N6").Value = 1.321648
$ws.Range("O6").Value = 0.00232403844212461
$ws.Range("P6").Value = 0.00232403844212461
$ws.Range("Q6").Value = 8.921533710880002
$ws.Range("R6").Value = 80.29380339792002
$ws.Range("S6").Value = 0.001800704642093025
$ws.Range("T6").Value = 0.001800704642093024
$ws.Range("G7").Value = 5.644020333333333
$ws.Range("I7").Value = 0.2159448115817524
$ws.Range("J7").Value = 0.2159448115817524
$ws.Range("M7").Value = 186.4969126666666
$ws.Range("N7").Value = 559.490738
$ws.Range("O7").Value = 0.9838307803020685
$ws.Range("P7").Value = 0.9838307803020684
$ws.Range("Q7").Value = 1052.592367194557
$ws.Range("R7").Value = 9473.331304751016
$ws.Range("S7").Value = 0.2124531524806586
$ws.Range("T7").Value = 0.2124531524806585
$ws.Range("G8").Value = 5.644020333333333
$ws.Range("I8").Value = 0.2159448115817524
$ws.Range("J8").Value = 0.2159448115817524
$ws.Range("O8").Value = 0.004224423017480373
$ws.Range("P8").Value = 0.004224423017480372
$ws.Range("Q8").Value = 4.519675042729999
$ws.Range("R8").Value = 40.67707538456999
$ws.Range("S8").Value = 0.000912242232551417
$ws.Range("T8").Value = 0.0009122422325514167
$ws.Range("G9").Value = 5.644020333333333
$ws.Range("I9").Value = 0.2159448115817524
$ws.Range("J9").Value = 0.2159448115817524
$ws.Range("M9").Value = 1.711075333333334
$ws.Range("N9").Value = 5.133226000000001
$ws.Range("O9").Value = 0.009026468890440984
$ws.Range("P9").Value = 0.009026468890440982
$ws.Range("Q9").Value = 9.657343973198444
$ws.Range("R9").Value = 86.916095758786
$ws.Range("S9").Value = 0.001949219123794828
$ws.Range("T9").Value = 0.001949219123794827
$ws.Range("G10").Value = 5.644020333333333
$ws.Range("I10").Value = 0.2159448115817524
$ws.Range("J10").Value = 0.2159448115817524
$ws.Range("M10").Value = 0.1126546666666667
$ws.Range("N10").Value = 0.337964
$ws.Range("O10").Value = 0.00059428934788552
$ws.Range("P10").Value = 0.0005942893478855199
$ws.Range("Q10").Value = 0.6358252293115555
$ws.Range("R10").Value = 5.722427063803999
$ws.Range("S10").Value = 0.0001283337012541811
$ws.Range("T10").Value = 0.0001283337012541811
$ws.Range("G11").Value = 5.644020333333333
$ws.Range("I11").Value = 0.2159448115817524
$ws.Range("J11").Value = 0.2159448115817524
$ws.Range("M11").Value = 0.4405493333333334
$ws.Range("N11").Value = 1.321648
$ws.Range("O11").Value = 0.00232403844212461
$ws.Range("P11").Value = 0.00232403844212461
$ws.Range("Q11").Value = 2.486469395169778
$ws.Range("R11").Value = 22.378224556528
$ws.Range("S11").Value = 0.0005018640434933483
$ws.Range("T11").Value = 0.0005018640434933482
$ws.Range("A12").Value = 'MuSCs'
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.216663
$ws.Range("H12").Value = 0.649989
$ws.Range("I12").Value = 0.008289702720490535
$ws.Range("J12").Value = 0.008289702720490533
$ws.Range("M12").Value = 186.4969126666666
$ws.Range("N12").Value = 559.490738
$ws.Range("O12").Value = 0.9838307803020685
$ws.Range("P12").Value = 0.9838307803020684
$ws.Range("Q12").Value = 40.406980589098
$ws.Range("R12").Value = 363.662825301882
$ws.Range("S12").Value = 0.008155664695972382
$ws.Range("T12").Value = 0.00815566469597238
$ws.Range("A13").Value = 'MuSCs'
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.216663
$ws.Range("H13").Value = 0.649989
$ws.Range("I13").Value = 0.008289702720490535
$ws.Range("J13").Value = 0.008289702720490533
$ws.Range("O13").Value = 0.004224423017480373
$ws.Range("P13").Value = 0.004224423017480372
$ws.Range("Q13").Value = 0.17350156377
$ws.Range("R13").Value = 1.56151407393
$ws.Range("S13").Value = 0.00003501921098050988
$ws.Range("T13").Value = 0.00003501921098050987
$ws.Range("A14").Value = 'MuSCs'
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.216663
$ws.Range("H14").Value = 0.649989
$ws.Range("I14").Value = 0.008289702720490535
$ws.Range("J14").Value = 0.008289702720490533
$ws.Range("M14").Value = 1.711075333333334
$ws.Range("N14").Value = 5.133226000000001
$ws.Range("O14").Value = 0.009026468890440984
$ws.Range("P14").Value = 0.009026468890440982
$ws.Range("Q14").Value = 0.3707267149460001
$ws.Range("R14").Value = 3.336540434514001
$ws.Range("S14").Value = 0.0000748267437175118
$ws.Range("T14").Value = 0.00007482674371751177
$ws.Range("A15").Value = 'MuSCs'
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.216663
$ws.Range("H15").Value = 0.649989
$ws.Range("I15").Value = 0.008289702720490535
$ws.Range("J15").Value = 0.008289702720490533
$ws.Range("M15").Value = 0.1126546666666667
$ws.Range("N15").Value = 0.337964
$ws.Range("O15").Value = 0.00059428934788552
$ws.Range("P15").Value = 0.0005942893478855199
$ws.Range("Q15").Value = 0.024408098044
$ws.Range("R15").Value = 0.219672882396
$ws.Range("S15").Value = 0.000004926482023925141
$ws.Range("T15").Value = 0.000004926482023925139
$ws.Range("A16").Value = 'MuSCs'
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.216663
$ws.Range("H16").Value = 0.649989
$ws.Range("I16").Value = 0.008289702720490535
$ws.Range("J16").Value = 0.008289702720490533
$ws.Range("M16").Value = 0.4405493333333334
$ws.Range("N16").Value = 1.321648
$ws.Range("O16").Value = 0.00232403844212461
$ws.Range("P16").Value = 0.00232403844212461
$ws.Range("Q16").Value = 0.09545074020800003
$ws.Range("R16").Value = 0.8590566618720001
$ws.Range("S16").Value = 0.00001926558779620497
$ws.Range("T16").Value = 0.00001926558779620496
$ws.Range("A17").Value = 'Resolving-Mac'
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.02478833333333333
$ws.Range("H17").Value = 0.074365
$ws.Range("I17").Value = 0.0009484218083833396
$ws.Range("J17").Value = 0.0009484218083833395
$ws.Range("M17").Value = 186.4969126666666
$ws.Range("N17").Value = 559.490738
$ws.Range("O17").Value = 0.9838307803020685
$ws.Range("P17").Value = 0.9838307803020684
$ws.Range("Q17").Value = 4.622947636818888
$ws.Range("R17").Value = 41.60652873137
$ws.Range("S17").Value = 0.0009330865677972799
$ws.Range("T17").Value = 0.0009330865677972796
$ws.Range("A18").Value = 'Resolving-Mac'
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.02478833333333333
$ws.Range("H18").Value = 0.074365
$ws.Range("I18").Value = 0.0009484218083833396
$ws.Range("J18").Value = 0.0009484218083833395
$ws.Range("O18").Value = 0.004224423017480373
$ws.Range("P18").Value = 0.004224423017480372
$ws.Range("Q18").Value = 0.01985024945
$ws.Range("R18").Value = 0.17865224505
$ws.Range("S18").Value = 0.000004006534917614939
$ws.Range("T18").Value = 0.000004006534917614938
$ws.Range("A19").Value = 'Resolving-Mac'
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.02478833333333333
$ws.Range("H19").Value = 0.074365
$ws.Range("I19").Value = 0.0009484218083833396
$ws.Range("J19").Value = 0.0009484218083833395
$ws.Range("M19").Value = 1.711075333333334
$ws.Range("N19").Value = 5.133226000000001
$ws.Range("O19").Value = 0.009026468890440984
$ws.Range("P19").Value = 0.009026468890440982
$ws.Range("Q19").Value = 0.04241470572111111
$ws.Range("R19").Value = 0.3817323514900001
$ws.Range("S19").Value = 0.000008560899948387995
$ws.Range("T19").Value = 0.000008560899948387992
$ws.Range("A20").Value = 'Resolving-Mac'
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.02478833333333333
$ws.Range("H20").Value = 0.074365
$ws.Range("I20").Value = 0.0009484218083833396
$ws.Range("J20").Value = 0.0009484218083833395
$ws.Range("M20").Value = 0.1126546666666667
$ws.Range("N20").Value = 0.337964
$ws.Range("O20").Value = 0.00059428934788552
$ws.Range("P20").Value = 0.0005942893478855199
$ws.Range("Q20").Value = 0.002792521428888889
$ws.Range("R20").Value = 0.02513269286
$ws.Range("S20").Value = 0.0000005636369780245405
$ws.Range("T20").Value = 0.0000005636369780245403
$ws.Range("A21").Value = 'Resolving-Mac'
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.02478833333333333
$ws.Range("H21").Value = 0.074365
$ws.Range("I21").Value = 0.0009484218083833396
$ws.Range("J21").Value = 0.0009484218083833395
$ws.Range("M21").Value = 0.4405493333333334
$ws.Range("N21").Value = 1.321648
$ws.Range("O21").Value = 0.00232403844212461
$ws.Range("P21").Value = 0.00232403844212461
$ws.Range("Q21").Value = 0.01092048372444445
$ws.Range("R21").Value = 0.09828435352
$ws.Range("S21").Value = 0.000002204168742032222
$ws.Range("T21").Value = 0.000002204168742032222
